# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holding detail) right before the
#    "总计" (totals) sheet, matching the existing per-quarter sheet layout.
# 2. Insert a new summary row at the top of the "总计" sheet's data for the
#    2022-Q1 quarter and renumber the existing index column.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" sheet, positioned immediately before "总计"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# NOTE: Worksheets.Add(Before) reseats the "Before" reference it was handed
# onto the freshly created sheet, so re-resolve "总计" by name before using
# it again below.
$total = $wb.Worksheets.Item("总计")

# Pull over the header-row formatting (bold / border / centered) and the
# index-column formatting from an existing quarter sheet so the new sheet
# matches the workbook's established look.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

# Fund codes (B) and the numeric-looking text metrics (D:G) must stay text so
# leading zeros / trailing zeros in codes and percentages survive.
$newSheet.Range("B2:G6").NumberFormat = "@"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4

$newSheet.Range("B2").Value = "161724"
$newSheet.Range("C2").Value = "招商中证煤炭等权指数（LOF）"
$newSheet.Range("D2").Value = "21.53"
$newSheet.Range("E2").Value = "94.64"
$newSheet.Range("F2").Value = "3.20"
$newSheet.Range("G2").Value = "0.6890"
$newSheet.Range("H2").Value = 10

$newSheet.Range("B3").Value = "217024"
$newSheet.Range("C3").Value = "招商安盈债券"
$newSheet.Range("D3").Value = "35.05"
$newSheet.Range("E3").Value = "20.20"
$newSheet.Range("F3").Value = "1.40"
$newSheet.Range("G3").Value = "0.4907"
$newSheet.Range("H3").Value = 3

$newSheet.Range("B4").Value = "014887"
$newSheet.Range("C4").Value = "招商安福1年定期开放债券"
$newSheet.Range("D4").Value = "17.22"
$newSheet.Range("E4").Value = "27.65"
$newSheet.Range("F4").Value = "1.37"
$newSheet.Range("G4").Value = "0.2359"
$newSheet.Range("H4").Value = 7

$newSheet.Range("B5").Value = "001219"
$newSheet.Range("C5").Value = "上投摩根动态多因子策略混合"
$newSheet.Range("D5").Value = "1.16"
$newSheet.Range("E5").Value = "92.44"
$newSheet.Range("F5").Value = "3.62"
$newSheet.Range("G5").Value = "0.0420"
$newSheet.Range("H5").Value = 2

$newSheet.Range("B6").Value = "519097"
$newSheet.Range("C6").Value = "新华中小市值优选混合"
$newSheet.Range("D6").Value = "0.75"
$newSheet.Range("E6").Value = "62.70"
$newSheet.Range("F6").Value = "2.24"
$newSheet.Range("G6").Value = "0.0168"
$newSheet.Range("H6").Value = 10

# ---------------------------------------------------------------------------
# 2. Prepend a 2022-Q1 row to the "总计" summary sheet
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Restore the index-column formatting on the newly inserted row, then clear
# the stray formatting Insert() carried into the new data cells.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 1.47

# Renumber the existing rows' index column (they shifted down by one row).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
